$d = $word.ActiveDocument

# 1. "τον teammate του να περάσει ούτε" -> "τον teammate του να τον περάσει ούτε"
$d.Content.Find.Execute("teammate του να περάσει ούτε", $true, $false, $false, $false, $false,
                         $true, 1, $false, "teammate του να τον περάσει ούτε", 2) | Out-Null

# 2. "Όλο γύρω πάνε!" -> "Όλο γύρω γύρω πάνε!"
$d.Content.Find.Execute("Όλο γύρω πάνε!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Όλο γύρω γύρω πάνε!", 2) | Out-Null

# 3. ", για να μην χάσει τη γραμμή του" -> ", για να μην χάσει τη γραμμή του ερχομενος με 360 km/h !!!"
$d.Content.Find.Execute(", για να μην χάσει τη γραμμή του", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", για να μην χάσει τη γραμμή του ερχομενος με 360 km/h !!!", 2) | Out-Null

# 4. "που έφαγε curb λες και ήταν μπισκότο" -> "που έφαγε το σασι στα curb λες και εξαρταται η ζωη του οδηγου απο αυτο."
$d.Content.Find.Execute("που έφαγε curb λες και ήταν μπισκότο", $true, $false, $false, $false, $false,
                         $true, 1, $false, "που έφαγε το σασι στα curb λες και εξαρταται η ζωη του οδηγου απο αυτο.", 2) | Out-Null

# 5. "F1 είναι το Game of Thrones" -> "Η F1 είναι το Game of Thrones"
$d.Content.Find.Execute("F1 είναι το Game of Thrones", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Η F1 είναι το Game of Thrones", 2) | Out-Null

# 6. "Δεν θες άλλον έναν που:" -> two runs "Σιγουρα" + " θελεις άλλον έναν που:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Δεν θες άλλον έναν που:") {
        $para.Range.Text = "Σιγουρα θελεις άλλον έναν που:"
        $splitPoint = $d.Range($para.Range.Start, $para.Range.Start + 7)
        $splitPoint.Bold = 1
        $splitPoint.Bold = 0
        break
    }
}

# 7. "Να πετάει memes για τον Toto Wolff στα διαλείμματα" -> append " στην ομαδικη"
$d.Content.Find.Execute("Toto Wolff στα διαλείμματα", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Toto Wolff στα διαλείμματα στην ομαδικη", 2) | Out-Null

# 8. Split paragraph: move " Το υπόλοιπο 30% ..." sentence (after the <w:br/>) into a brand-new paragraph
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $full = $para.Range
    $txt = $full.Text
    $breakIdx = $txt.IndexOf([char]11)
    if ($breakIdx -ge 0 -and $txt.IndexOf("Το υπόλοιπο 30%") -gt $breakIdx) {
        $start = $full.Start
        $delStart = $start + $breakIdx + 1
        $delEnd = $start + $txt.Length - 1
        $delRange = $d.Range($delStart, $delEnd)
        $movedText = $delRange.Text.TrimStart()
        $delRange.Text = ""
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = $movedText
        break
    }
}
